# Add a new row of data (row 6) to Sheet1 to confirm the correct file is
# being read, matching the existing ID / First_Name / Last_Name / Position
# table layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Employee"
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = "Just proving we are reading the correct documents"

# Mirror the cursor ending up one row below/after the newly entered data,
# as happens after typing values across a row in Excel.
$ws.Range("D7").Select()
